# Fix latency units in report sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to the Read Latency columns (I, J, K) for data rows 3-23
for ($row = 3; $row -le 23; $row++) {
    $cellI = $ws.Range("I$row")
    $cellI.Value = "$($cellI.Value()) msec"

    $cellJ = $ws.Range("J$row")
    $cellJ.Value = "$($cellJ.Value()) msec"

    $cellK = $ws.Range("K$row")
    $cellK.Value = "$($cellK.Value()) msec"
}
